# Update "想去人数" (wanted-to-go count) figures in column F
# for the "展览" (Exhibition) and "全部类型" (All types) sheets.
# These two sheets mirror the same underlying event data, so the same
# set of row-level updates is applied to both (rows differ by one
# because "全部类型" contains an extra event not present on "展览").

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# Row => New value, for sheet "展览"
$exhibitionUpdates = @{
    2  = 1876
    3  = 492
    6  = 2590
    8  = 92
    10 = 1535
    11 = 531
    13 = 334
    14 = 231
    21 = 177
    22 = 60
    23 = 1665
    24 = 29
    25 = 407
    27 = 207
    28 = 302
    29 = 420
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row => New value, for sheet "全部类型"
$allTypesUpdates = @{
    2  = 1876
    4  = 492
    7  = 2590
    9  = 92
    11 = 1535
    12 = 531
    14 = 334
    15 = 231
    22 = 177
    23 = 60
    24 = 1665
    25 = 29
    26 = 407
    28 = 207
    29 = 302
    30 = 420
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
